$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are not
# reinterpreted as numbers (they must remain text, matching the source data).
$ws.Range("D2:D51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '27.697.11'
$ws.Range("D3").Value = '1.755.48'
$ws.Range("D5").Value = '324.41'
$ws.Range("D6").Value = '0.9986'
$ws.Range("D7").Value = '0.4282'
$ws.Range("D8").Value = '0.3641'
$ws.Range("D9").Value = '45.35'
$ws.Range("D10").Value = '0.07493'
$ws.Range("D11").Value = '1.123'
$ws.Range("D12").Value = '0.9993'
$ws.Range("D13").Value = '21.66'
$ws.Range("D14").Value = '6.151'
$ws.Range("D15").Value = '7.262'
$ws.Range("D16").Value = '1.748.31'
$ws.Range("D17").Value = '0.00001072'
$ws.Range("D18").Value = '87.72'
$ws.Range("D19").Value = '0.06205'
$ws.Range("D20").Value = '0.9989'
$ws.Range("D21").Value = '17.07'
$ws.Range("D22").Value = '6.159'
$ws.Range("D23").Value = '0.5263'
$ws.Range("D24").Value = '27.699.80'
$ws.Range("D26").Value = '2.336'
$ws.Range("D28").Value = '152.65'
$ws.Range("D29").Value = '2.375'
$ws.Range("D30").Value = '1.946.97'
$ws.Range("D33").Value = '5.757'
$ws.Range("D34").Value = '0.09155'
$ws.Range("D35").Value = '3.653'
$ws.Range("D36").Value = '12.78'
$ws.Range("D37").Value = '0.02319'
$ws.Range("D38").Value = '0.2154'
$ws.Range("D39").Value = '5.128'
$ws.Range("D40").Value = '0.6500'
$ws.Range("D41").Value = '0.06121'
$ws.Range("D42").Value = '1.199'
$ws.Range("D43").Value = '1.423'
$ws.Range("D44").Value = '7.981'
$ws.Range("D45").Value = '0.9982'
$ws.Range("D46").Value = '13.85'
$ws.Range("D47").Value = '0.5949'
$ws.Range("D48").Value = '3.738'
$ws.Range("D49").Value = '126.32'
$ws.Range("D50").Value = '1.974'
$ws.Range("D51").Value = '0.06903'

# Restore default (Normal) style for column D so no extra number-format
# styling is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"

# --- Column B (Coin) updates: rows 28/29 swapped (LidoDAOToken <-> Monero) ---
$ws.Range("B28").Value = 'Monero'
$ws.Range("B29").Value = 'LidoDAOToken'

# --- Column C (Link) updates: rows 28/29 swapped ---
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  -2.05%  '
$ws.Range("E3").Value = '  -2.82%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("E5").Value = '  -4.12%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("E7").Value = '  -7.83%  '
$ws.Range("E8").Value = '  -4.81%  '
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  -3.25%  '
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("E13").Value = '  -3.88%  '
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("E15").Value = '  -3.37%  '
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("E18").Value = '  +7.47%  '
$ws.Range("E19").Value = '  -7.85%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("E22").Value = '  -4.32%  '
$ws.Range("E23").Value = '  -4.72%  '
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("E26").Value = '  -3.87%  '
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("E32").Value = '  -4.30%  '
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("E34").Value = '  -4.98%  '
$ws.Range("E35").Value = '  -9.55%  '
$ws.Range("E36").Value = '  +5.43%  '
$ws.Range("E37").Value = '  -1.71%  '
$ws.Range("E38").Value = '  -8.46%  '
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("E41").Value = '  -3.88%  '
$ws.Range("E42").Value = '  -3.80%  '
$ws.Range("E43").Value = '  -4.66%  '
$ws.Range("E44").Value = '  -4.61%  '
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("E46").Value = '  -3.16%  '
$ws.Range("E47").Value = '  -3.28%  '
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("E49").Value = '  -4.13%  '
$ws.Range("E50").Value = '  -3.93%  '
$ws.Range("E51").Value = '  -3.67%  '
